$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column F ("język") is being removed entirely; delete the whole column
# so everything to its right (grupa, nr tel, email, notatka rekrutacyjna)
# shifts one column to the left.
$ws.Range("F1").EntireColumn.Delete()

# Match the selection left behind by the author's edit.
$ws.Range("F1").Select()
